# Apply updated cryptocurrency market data to Sheet1
# (values are stored as text in the source sheet, so numeric-looking
#  price strings are force-formatted as Text before assignment to avoid
#  Excel auto-converting them to numbers and losing trailing zeros / exact text)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.652.01"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "1.633.90"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.90"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  +3.48%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0623"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.18"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("E11").Value = "  +3.35%  "
$ws.Range("D12").Value = "1.860.60"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "1.643.46"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D16").Value = "26.638.53"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.35"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.97"
$ws.Range("E19").Value = "  +8.56%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.29"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.20"
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.43"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.56"
$ws.Range("E25").Value = "  +2.62%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.90"
$ws.Range("E28").Value = "  +5.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.47"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.52"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "1.211.20"
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0173"
$ws.Range("E37").Value = "  +5.59%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.503"
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").Value = "1.769.65"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.86"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.69"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0512"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.65"
$ws.Range("E49").Value = "  +5.06%  "
$ws.Range("E50").Value = "  +0.35%  "
